# Insert a new worksheet "具有相當價值之財產" immediately before the "保險" sheet,
# carrying one row describing an antique (朱銘雕刻) owned by 陳明文 worth 500000.
$wb = $excel.ActiveWorkbook

$insuranceSheet = $wb.Worksheets.Item("保險")

# Create the new sheet positioned right before "保險" so the tab order becomes:
# 土地, 建物, 存款, 股票, 具有相當價值之財產, 保險
$newSheet = $wb.Worksheets.Add($insuranceSheet)
$newSheet.Name = "具有相當價值之財產"

# ---- Header row (row 1): bold, centered, top-aligned, thin border all round
# (matches the header style used on every other sheet in the workbook) ----
$headerRange = $newSheet.Range("B1:L1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

$newSheet.Range("B1").Value = "name"
$newSheet.Range("C1").Value = "quantity"
$newSheet.Range("D1").Value = "owner"
$newSheet.Range("E1").Value = "total"
$newSheet.Range("F1").Value = "property_category"
$newSheet.Range("G1").Value = "category"
$newSheet.Range("H1").Value = "date"
$newSheet.Range("I1").Value = "legislator_name"
$newSheet.Range("J1").Value = "legislator_id"
$newSheet.Range("K1").Value = "source_file"
$newSheet.Range("L1").Value = "index"

# ---- Data row (row 2): plain font, thin border all round ----
$dataRange = $newSheet.Range("A2:L2")
$dataRange.Borders.LineStyle = 1

# Column A (the numeric index) follows the same convention used on every
# other sheet in this workbook: it reuses the bold/centered header style
# even though it sits in a data row.
$newSheet.Range("A2").Font.Bold = $true
$newSheet.Range("A2").HorizontalAlignment = -4108   # xlCenter
$newSheet.Range("A2").VerticalAlignment = -4160     # xlTop

$newSheet.Range("A2").Value = 142
$newSheet.Range("B2").Value = "朱銘雕刻"
$newSheet.Range("C2").Value = 1
$newSheet.Range("D2").Value = "陳明文"
$newSheet.Range("E2").Value = 500000
$newSheet.Range("F2").Value = "antique"
$newSheet.Range("G2").Value = "normal"

# The "date" column stores a literal text value ("2013-11-20"), not an Excel
# date serial, so force text formatting before assigning it.
$newSheet.Range("H2").NumberFormat = "@"
$newSheet.Range("H2").Value = "2013-11-20"

$newSheet.Range("I2").Value = "陳明文"
$newSheet.Range("J2").Value = 828
$newSheet.Range("K2").Value = "tmp581f1"
$newSheet.Range("L2").Value = 142

# Restore the workbook's original active sheet ("土地") so adding/naming the
# new sheet doesn't leave a different tab focused.
$wb.Worksheets.Item("土地").Activate()

